$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-12 Tuesday", "2024-11-13 Wednesday"),
    @("84×97=8148", "56×24=1344"),
    @("98×57=5586", "47×26=1222"),
    @("24×27=648", "18×40=720"),
    @("91×45=4095", "72×83=5976"),
    @("90×87=7830", "88×41=3608"),
    @("16×40=640", "97×31=3007"),
    @("98×92=9016", "46×24=1104"),
    @("59×57=3363", "58×59=3422"),
    @("30×65=1950", "82×17=1394"),
    @("80×47=3760", "42×43=1806"),
    @("18×52=936", "20×34=680"),
    @("37×76=2812", "87×56=4872"),
    @("32×91=2912", "95×83=7885"),
    @("74×42=3108", "78×95=7410"),
    @("82×84=6888", "85×14=1190"),
    @("39×19=741", "30×11=330"),
    @("84×98=8232", "72×96=6912"),
    @("55×52=2860", "40×98=3920"),
    @("59×13=767", "97×34=3298"),
    @("11×65=715", "32×13=416"),
    @("16×50=800", "12×94=1128"),
    @("63×61=3843", "67×84=5628"),
    @("70×72=5040", "65×25=1625"),
    @("76×95=7220", "28×55=1540"),
    @("87×34=2958", "88×27=2376")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
